$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to metrics_sim_with_priors.json - correcting relevance markers
$ws.Range("C3").Value = 0.3820224719101123
$ws.Range("D3").Value = 0.6629213483146067
$ws.Range("E3").Value = 0.9101123595505618
$ws.Range("H3").Value = 0.3889655172413793
$ws.Range("I3").Value = 0.1690339905307045
$ws.Range("J3").Value = 0.2808988764044944
$ws.Range("K3").Value = 152.5056179775281

$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 24
$ws.Range("S3").Value = 64
$ws.Range("T3").Value = 161
$ws.Range("U3").Value = 268
$ws.Range("V3").Value = 629
$ws.Range("W3").Value = 612
$ws.Range("X3").Value = 572
$ws.Range("Y3").Value = 475
$ws.Range("Z3").Value = 368

$ws.Range("AF3").Value = 0.988994
$ws.Range("AG3").Value = 0.962264
$ws.Range("AH3").Value = 0.899371
$ws.Range("AI3").Value = 0.746855
$ws.Range("AJ3").Value = 0.578616
